$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header text in F1 from "Draft" to "Drafting of manuscript"
$ws.Range("F1").Value = "Drafting of manuscript"

# Reflect the new active cell selection (user clicked F2 after editing)
$ws.Range("F2").Select()
